# Applies the periodic cryptos-list price/volume refresh described by the
# commit "Updated cryptos list ... with GitHub Actions".
#
# Numeric-looking text (e.g. "210.99") is written with a leading apostrophe,
# exactly as a user typing into Excel would, so the cell keeps storing plain
# text instead of Excel auto-converting it to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '27.324.54'
$ws.Range('E2').Value = '  +1.45%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '1.570.51'
$ws.Range('E3').Value = '  +0.70%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.23%  '

# Row 5: BNB
$ws.Range('D5').Value = '''210.99'
$ws.Range('E5').Value = '  +1.87%  '

# Row 6: XRP
$ws.Range('D6').Value = '''0.493'
$ws.Range('E6').Value = '  +0.79%  '

# Row 7: USDC
$ws.Range('E7').Value = '  +0.11%  '

# Row 8: Solana
$ws.Range('D8').Value = '''22.19'
$ws.Range('E8').Value = '  +0.68%  '

# Row 9: Cardano
$ws.Range('D9').Value = '''0.250'
$ws.Range('E9').Value = '  +0.71%  '

# Row 10: Dogecoin
$ws.Range('E10').Value = '  +0.13%  '

# Row 11: TRON
$ws.Range('E11').Value = '  +1.45%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Range('D12').Value = '1.794.00'
$ws.Range('E12').Value = '  +0.66%  '

# Row 13: WrappedEther
$ws.Range('D13').Value = '1.561.59'
$ws.Range('E13').Value = '  +0.07%  '

# Row 14: Polkadot
$ws.Range('E14').Value = '  +0.93%  '

# Row 15: Polygon
$ws.Range('E15').Value = '  +0.43%  '

# Row 16: WrappedBTC
$ws.Range('D16').Value = '27.263.48'
$ws.Range('E16').Value = '  +1.18%  '

# Row 17: Litecoin
$ws.Range('D17').Value = '''62.28'
$ws.Range('E17').Value = '  +0.41%  '

# Row 18: Chainlink
$ws.Range('D18').Value = '''7.55'
$ws.Range('E18').Value = '  +2.64%  '

# Row 19: BitcoinCash
$ws.Range('D19').Value = '''217.45'
$ws.Range('E19').Value = '  +0.48%  '

# Row 20: ShibaInu
$ws.Range('E20').Value = '  -0.07%  '

# Row 21: Dai
$ws.Range('E21').Value = '  +0.16%  '

# Row 22: Uniswap
$ws.Range('E22').Value = '  +1.52%  '

# Row 23: Avalanche
$ws.Range('D23').Value = '''9.25'
$ws.Range('E23').Value = '  +0.45%  '

# Row 24: Toncoin
$ws.Range('E24').Value = '  +0.33%  '

# Row 25: Monero
$ws.Range('D25').Value = '''153.74'

# Row 26: Cosmos
$ws.Range('E26').Value = '  +0.81%  '

# Row 27: EthereumClassic
$ws.Range('D27').Value = '''15.10'
$ws.Range('E27').Value = '  +0.47%  '

# Row 28: Stellar
$ws.Range('E28').Value = '  +2.25%  '

# Row 29: BinanceUSD
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  +0.19%  '

# Row 30: PancakeSwap
$ws.Range('E30').Value = '  +2.94%  '

# Row 31: Hedera
$ws.Range('E31').Value = '  +0.58%  '

# Row 32: Filecoin
$ws.Range('E32').Value = '  +0.53%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range('E33').Value = '  +1.75%  '

# Row 34: Maker
$ws.Range('D34').Value = '1.445.82'
$ws.Range('E34').Value = '  +2.01%  '

# Row 35: TrustWalletToken
$ws.Range('D35').Value = '''1.11'
$ws.Range('E35').Value = '  +4.32%  '

# Row 36: LidoDAOToken
$ws.Range('D36').Value = '''1.61'
$ws.Range('E36').Value = '  +0.18%  '

# Row 37: HuobiToken
$ws.Range('E37').Value = '  +0.46%  '

# Row 38: VeChain
$ws.Range('E38').Value = '  +1.04%  '

# Row 39: ImmutableX
$ws.Range('E39').Value = '  +0.63%  '

# Row 40: FraxShare
$ws.Range('D40').Value = '''5.87'
$ws.Range('E40').Value = '  +2.52%  '

# Row 41: MXToken
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = '''0.813'
$ws.Range('E41').Value = '  +0.71%  '

# Row 42: ARBITRUM
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = '''1.00'
$ws.Range('E42').Value = '  +0.12%  '

# Row 43: PaxDollar
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = '''2.34'
$ws.Range('E43').Value = '  +0.74%  '

# Row 44: WEMIXToken
$ws.Range('E44').Value = '  -0.67%  '

# Row 45: Aave
$ws.Range('E45').Value = '  +0.01%  '

# Row 46: RenderToken
$ws.Range('D46').Value = '''1.74'
$ws.Range('E46').Value = '  -0.68%  '

# Row 47: RocketPoolETH
$ws.Range('D47').Value = '1.706.12'
$ws.Range('E47').Value = '  +0.57%  '

# Row 48: Quant
$ws.Range('D48').Value = '''86.23'
$ws.Range('E48').Value = '  -1.28%  '

# Row 49: Cronos
$ws.Range('D49').Value = '''0.0526'
$ws.Range('E49').Value = '  +1.14%  '

# Row 50: BabyDogeCoin
$ws.Range('E50').Value = '  +1.05%  '

# Row 51: Algorand
$ws.Range('E51').Value = '  +0.57%  '
